$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "298.63"
Set-TextValue $ws.Range("E2") "-2.52%"
Set-TextValue $ws.Range("G2") "6"

Set-TextValue $ws.Range("D3") "31.74"
Set-TextValue $ws.Range("E3") "-1.61%"
Set-TextValue $ws.Range("G3") "6"

Set-TextValue $ws.Range("D4") "5.102"
Set-TextValue $ws.Range("E4") "-4.52%"
Set-TextValue $ws.Range("G4") "6"

Set-TextValue $ws.Range("D5") "0.07540"
Set-TextValue $ws.Range("E5") "1.60%"
Set-TextValue $ws.Range("G5") "6"

Set-TextValue $ws.Range("D6") "7.752"
Set-TextValue $ws.Range("E6") "-0.37%"
Set-TextValue $ws.Range("G6") "6"

Set-TextValue $ws.Range("D7") "1.765"
Set-TextValue $ws.Range("E7") "11.33%"
Set-TextValue $ws.Range("G7") "6"

Set-TextValue $ws.Range("D8") "3.796"
Set-TextValue $ws.Range("E8") "2.83%"
Set-TextValue $ws.Range("G8") "6"

Set-TextValue $ws.Range("D9") "0.9270"
Set-TextValue $ws.Range("E9") "1.59%"
Set-TextValue $ws.Range("G9") "6"

Set-TextValue $ws.Range("D10") "0.1706"
Set-TextValue $ws.Range("E10") "1.82%"
Set-TextValue $ws.Range("G10") "6"

Set-TextValue $ws.Range("D11") "0.07478"
Set-TextValue $ws.Range("E11") "-2.58%"
Set-TextValue $ws.Range("G11") "6"

Set-TextValue $ws.Range("D12") "0.07968"
Set-TextValue $ws.Range("E12") "-1.03%"
Set-TextValue $ws.Range("G12") "6"

Set-TextValue $ws.Range("D13") "0.03060"
Set-TextValue $ws.Range("E13") "1.41%"
Set-TextValue $ws.Range("G13") "6"

Set-TextValue $ws.Range("D14") "0.09901"
Set-TextValue $ws.Range("E14") "0.62%"
Set-TextValue $ws.Range("G14") "6"

Set-TextValue $ws.Range("D15") "0.001508"
Set-TextValue $ws.Range("E15") "-0.79%"
Set-TextValue $ws.Range("G15") "6"

Set-TextValue $ws.Range("D16") "0.006503"
Set-TextValue $ws.Range("E16") "5.94%"
Set-TextValue $ws.Range("G16") "6"

Set-TextValue $ws.Range("D17") "3.451"
Set-TextValue $ws.Range("G17") "6"

Set-TextValue $ws.Range("D18") "2.223"
Set-TextValue $ws.Range("E18") "-0.77%"
Set-TextValue $ws.Range("G18") "6"

Set-TextValue $ws.Range("E19") "0.49%"
Set-TextValue $ws.Range("G19") "6"

Set-TextValue $ws.Range("E20") "-0.72%"
Set-TextValue $ws.Range("G20") "6"

Set-TextValue $ws.Range("D21") "4.557"
Set-TextValue $ws.Range("E21") "8.51%"
Set-TextValue $ws.Range("G21") "6"

Set-TextValue $ws.Range("D22") "0.04650"
Set-TextValue $ws.Range("E22") "2.18%"
Set-TextValue $ws.Range("G22") "6"

Set-TextValue $ws.Range("E23") "-4.38%"
Set-TextValue $ws.Range("G23") "6"

Set-TextValue $ws.Range("D24") "0.001219"
Set-TextValue $ws.Range("E24") "0.37%"
Set-TextValue $ws.Range("G24") "6"

Set-TextValue $ws.Range("D25") "0.004422"
Set-TextValue $ws.Range("E25") "-1.69%"
Set-TextValue $ws.Range("G25") "6"

Set-TextValue $ws.Range("E26") "19.58%"
Set-TextValue $ws.Range("G26") "6"

Set-TextValue $ws.Range("D27") "0.0001858"
Set-TextValue $ws.Range("E27") "6.84%"
Set-TextValue $ws.Range("G27") "6"

Set-TextValue $ws.Range("G28") "6"

Set-TextValue $ws.Range("G29") "6"

Set-TextValue $ws.Range("G30") "6"

Set-TextValue $ws.Range("G31") "6"

Set-TextValue $ws.Range("G32") "6"

Set-TextValue $ws.Range("G33") "6"

Set-TextValue $ws.Range("G34") "6"

Set-TextValue $ws.Range("G35") "6"

Set-TextValue $ws.Range("G36") "6"

Set-TextValue $ws.Range("G37") "6"

Set-TextValue $ws.Range("G38") "6"

Set-TextValue $ws.Range("D39") "0.01686"
Set-TextValue $ws.Range("E39") "-1.01%"
Set-TextValue $ws.Range("G39") "6"

Set-TextValue $ws.Range("D40") "0.04550"
Set-TextValue $ws.Range("E40") "0.86%"
Set-TextValue $ws.Range("G40") "6"

Set-TextValue $ws.Range("D41") "0.007044"
Set-TextValue $ws.Range("E41") "-1.44%"
Set-TextValue $ws.Range("G41") "6"

Set-TextValue $ws.Range("E42") "-2.42%"
Set-TextValue $ws.Range("G42") "6"

Set-TextValue $ws.Range("E43") "-8.92%"
Set-TextValue $ws.Range("G43") "6"

Set-TextValue $ws.Range("D44") "0.01282"
Set-TextValue $ws.Range("E44") "-6.01%"
Set-TextValue $ws.Range("G44") "6"

Set-TextValue $ws.Range("D45") "0.00006040"
Set-TextValue $ws.Range("E45") "-1.31%"
Set-TextValue $ws.Range("G45") "6"

Set-TextValue $ws.Range("D46") "0.7116"
Set-TextValue $ws.Range("E46") "-62.40%"
Set-TextValue $ws.Range("G46") "6"

Set-TextValue $ws.Range("E47") "-0.15%"
Set-TextValue $ws.Range("G47") "6"

Set-TextValue $ws.Range("G48") "6"

Set-TextValue $ws.Range("G49") "6"

Set-TextValue $ws.Range("G50") "6"

Set-TextValue $ws.Range("G51") "6"
